$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.776.96'
$ws.Range("E2").Value = '  -4.99%  '
$ws.Range("D3").Value = '3.504.72'
$ws.Range("E3").Value = '  -6.07%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.997'
$ws.Range("E4").Value = '  -0.29%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '564.63'
$ws.Range("E5").Value = '  -7.79%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '184.99'
$ws.Range("E6").Value = '  -4.08%  '
$ws.Range("D7").Value = '3.496.07'
$ws.Range("E7").Value = '  -6.18%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.603'
$ws.Range("E8").Value = '  -5.73%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.996'
$ws.Range("E9").Value = '  -0.48%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.654'
$ws.Range("E10").Value = '  -10.44%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.143'
$ws.Range("E11").Value = '  -11.75%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '52.21'
$ws.Range("E12").Value = '  -13.36%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000253'
$ws.Range("E13").Value = '  -13.16%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.61'
$ws.Range("E14").Value = '  -9.28%  '
$ws.Range("D15").Value = '4.046.12'
$ws.Range("E15").Value = '  -6.41%  '
$ws.Range("E16").Value = '  -1.68%  '
$ws.Range("D17").Value = '3.479.09'
$ws.Range("E17").Value = '  -6.69%  '
$ws.Range("D18").Value = '65.543.09'
$ws.Range("E18").Value = '  -5.10%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.85'
$ws.Range("E19").Value = '  -8.85%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.81'
$ws.Range("E20").Value = '  -9.04%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.04'
$ws.Range("E21").Value = '  -10.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '381.91'
$ws.Range("E22").Value = '  -7.67%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.13'
$ws.Range("E23").Value = '  -10.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.07'
$ws.Range("E24").Value = '  -6.85%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.81'
$ws.Range("E25").Value = '  -5.34%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.84'
$ws.Range("E26").Value = '  -8.22%  '
$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.07'
$ws.Range("E27").Value = '  -7.15%  '
$ws.Range("B28").Value = 'LEO'
$ws.Range("C28").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.99'
$ws.Range("E28").Value = '  -1.19%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.48'
$ws.Range("E29").Value = '  -8.73%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.70'
$ws.Range("E30").Value = '  -11.32%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '30.71'
$ws.Range("E31").Value = '  -7.10%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.28'
$ws.Range("E32").Value = '  -7.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '12.02'
$ws.Range("E33").Value = '  -6.29%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '608.29'
$ws.Range("E34").Value = '  -6.54%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '62.49'
$ws.Range("E35").Value = '  -6.86%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.111'
$ws.Range("E36").Value = '  -10.00%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '41.03'
$ws.Range("E37").Value = '  -10.80%  '
$ws.Range("E38").Value = '  +0.19%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.394'
$ws.Range("E39").Value = '  -6.16%  '
$ws.Range("B40").Value = 'PEPE'
$ws.Range("C40").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D40").Value = '0.0₃0721'
$ws.Range("E40").Value = '  -13.89%  '
$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.997'
$ws.Range("E41").Value = '  -0.36%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.129'
$ws.Range("E42").Value = '  -8.67%  '
$ws.Range("D43").Value = '2.986.08'
$ws.Range("E43").Value = '  +2.56%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.77'
$ws.Range("E44").Value = '  -10.17%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.47'
$ws.Range("E45").Value = '  -6.92%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.18'
$ws.Range("E46").Value = '  +2.43%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0399'
$ws.Range("E47").Value = '  -11.51%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.128'
$ws.Range("E48").Value = '  -8.97%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '136.53'
$ws.Range("E49").Value = '  -4.73%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.23'
$ws.Range("E50").Value = '  -11.16%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.69'
$ws.Range("E51").Value = '  -3.71%  '
